$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last data row (row 11) down into the new row 12
$ws.Range("A11:O11").Copy()
$ws.Range("A12:O12").PasteSpecial(-4122)

# Populate the new "TestCase11_validateDragAndDrop" record
$ws.Range("A12").Value = "TestCase11_validateDragAndDrop"
$ws.Range("B12").Value = "Validate Drag And Drop functionality."
$ws.Range("C12").Value = "kw_draganddrop"
$ws.Range("D12").Value = "N"
$ws.Range("J12").Value = "Others"
$ws.Range("K12").Value = "Drag and Drop"

# Match the original row's wrapped-text height
$ws.Rows.Item(12).RowHeight = 30

# Update the view to reflect scrolling down to the newly added row
$ws.Range("M15").Select() | Out-Null
